$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Sheet ALC, Row 54 (Leve Item ID 2174)
$ws_ALC.Range("H54").Value = 42679.8
$ws_ALC.Range("I54").Value = 33999
$ws_ALC.Range("J54").Value = 44850
$ws_ALC.Range("K54").Value = 33999
$ws_ALC.Range("L54").Value = 44850
$ws_ALC.Range("M54").Value = -33513
$ws_ALC.Range("N54").Value = -45822

# Sheet ALC, Row 62 (Leve Item ID 27781)
$ws_ALC.Range("H62").Value = 3149.8333
$ws_ALC.Range("I62").Value = 1556.8572
$ws_ALC.Range("J62").Value = 5380
$ws_ALC.Range("K62").Value = 1556.8572
$ws_ALC.Range("L62").Value = 5380
$ws_ALC.Range("M62").Value = -932.8571999999999
$ws_ALC.Range("N62").Value = -6628

# Sheet ALC, Row 65 (Leve Item ID 27781)
$ws_ALC.Range("H65").Value = 3149.8333
$ws_ALC.Range("I65").Value = 1556.8572
$ws_ALC.Range("J65").Value = 5380
$ws_ALC.Range("K65").Value = 7784.286
$ws_ALC.Range("L65").Value = 26900
$ws_ALC.Range("M65").Value = -4664.286
$ws_ALC.Range("N65").Value = -33140

# Sheet ALC, Row 70 (Leve Item ID 12604)
$ws_ALC.Range("H70").Value = 1328.0864
$ws_ALC.Range("I70").Value = 987.6923
$ws_ALC.Range("J70").Value = 2710.9375
$ws_ALC.Range("K70").Value = 2963.0769
$ws_ALC.Range("L70").Value = 8132.8125
$ws_ALC.Range("M70").Value = -2693.0769
$ws_ALC.Range("N70").Value = -8672.8125

# Sheet ALC, Row 73 (Leve Item ID 12604)
$ws_ALC.Range("H73").Value = 1328.0864
$ws_ALC.Range("I73").Value = 987.6923
$ws_ALC.Range("J73").Value = 2710.9375
$ws_ALC.Range("K73").Value = 2963.0769
$ws_ALC.Range("L73").Value = 8132.8125
$ws_ALC.Range("M73").Value = -2027.0769
$ws_ALC.Range("N73").Value = -10004.8125

# Sheet ALC, Row 112 (Leve Item ID 27960)
$ws_ALC.Range("H112").Value = 1625.5186
$ws_ALC.Range("J112").Value = 1648.9493
$ws_ALC.Range("L112").Value = 4946.8479
$ws_ALC.Range("N112").Value = -7162.8479

# Sheet ALC, Row 128 (Leve Item ID 34540)
$ws_ALC.Range("H128").Value = 42780
$ws_ALC.Range("J128").Value = 42780
$ws_ALC.Range("L128").Value = 42780
$ws_ALC.Range("N128").Value = -52740

# Sheet ALC, Row 129 (Leve Item ID 36115)
$ws_ALC.Range("H129").Value = 1031.1389
$ws_ALC.Range("J129").Value = 1175.8
$ws_ALC.Range("L129").Value = 3527.4
$ws_ALC.Range("N129").Value = -13527.4

# Sheet ALC, Row 130 (Leve Item ID 34691)
$ws_ALC.Range("H130").Value = 0
$ws_ALC.Range("J130").Value = 0
$ws_ALC.Range("L130").Value = 0
$ws_ALC.Range("N130").ClearContents()

# Sheet ALC, Row 133 (Leve Item ID 41856)
$ws_ALC.Range("H133").Value = 37025.715
$ws_ALC.Range("J133").Value = 37025.715
$ws_ALC.Range("L133").Value = 37025.715
$ws_ALC.Range("N133").Value = -47145.715

# Sheet ALC, Row 134 (Leve Item ID 41997)
$ws_ALC.Range("H134").Value = 42343.1
$ws_ALC.Range("J134").Value = 42343.1
$ws_ALC.Range("L134").Value = 42343.1
$ws_ALC.Range("N134").Value = -52483.1

# Sheet ALC, Row 136 (Leve Item ID 42164)
$ws_ALC.Range("H136").Value = 45803.793
$ws_ALC.Range("J136").Value = 45803.793
$ws_ALC.Range("L136").Value = 45803.793
$ws_ALC.Range("N136").Value = -56003.793

# Sheet ALC, Row 139 (Leve Item ID 42306)
$ws_ALC.Range("H139").Value = 41824.348
$ws_ALC.Range("J139").Value = 41824.348
$ws_ALC.Range("L139").Value = 41824.348
$ws_ALC.Range("N139").Value = -52104.348

# Sheet ALC, Row 140 (Leve Item ID 42459)
$ws_ALC.Range("H140").Value = 41175.758
$ws_ALC.Range("I140").Value = 39500
$ws_ALC.Range("J140").Value = 41283.87
$ws_ALC.Range("K140").Value = 39500
$ws_ALC.Range("L140").Value = 41283.87
$ws_ALC.Range("M140").Value = -34320
$ws_ALC.Range("N140").Value = -51643.87

# Sheet ARM, Row 32 (Leve Item ID 44147)
$ws_ARM.Range("H32").Value = 6725.396
$ws_ARM.Range("I32").Value = 4804.343
$ws_ARM.Range("K32").Value = 4804.343
$ws_ARM.Range("M32").Value = -4517.343

# Sheet ARM, Row 132 (Leve Item ID 43997)
$ws_ARM.Range("H132").Value = 2566.52
$ws_ARM.Range("I132").Value = 1976.2927
$ws_ARM.Range("J132").Value = 5255.3335
$ws_ARM.Range("K132").Value = 5928.8781
$ws_ARM.Range("L132").Value = 15766.0005
$ws_ARM.Range("M132").Value = -3398.8781
$ws_ARM.Range("N132").Value = -20826.0005

# Sheet BSM, Row 134 (Leve Item ID 43998)
$ws_BSM.Range("H134").Value = 2472.625
$ws_BSM.Range("I134").Value = 1480.8529
$ws_BSM.Range("K134").Value = 4442.5587
$ws_BSM.Range("M134").Value = -1907.5587

# Sheet CRP, Row 31 (Leve Item ID 44023)
$ws_CRP.Range("H31").Value = 11630945
$ws_CRP.Range("I31").Value = 1652.5172
$ws_CRP.Range("J31").Value = 35720190
$ws_CRP.Range("K31").Value = 1652.5172
$ws_CRP.Range("L31").Value = 35720190
$ws_CRP.Range("M31").Value = -1357.5172
$ws_CRP.Range("N31").Value = -35720780

# Sheet CRP, Row 34 (Leve Item ID 44023)
$ws_CRP.Range("H34").Value = 11630945
$ws_CRP.Range("I34").Value = 1652.5172
$ws_CRP.Range("J34").Value = 35720190
$ws_CRP.Range("K34").Value = 1652.5172
$ws_CRP.Range("L34").Value = 35720190
$ws_CRP.Range("M34").Value = -1450.5172
$ws_CRP.Range("N34").Value = -35720594

# Sheet CRP, Row 86 (Leve Item ID 12584)
$ws_CRP.Range("H86").Value = 2635
$ws_CRP.Range("I86").Value = 2589.25
$ws_CRP.Range("J86").Value = 2726.5
$ws_CRP.Range("K86").Value = 2589.25
$ws_CRP.Range("L86").Value = 2726.5
$ws_CRP.Range("M86").Value = -1466.25
$ws_CRP.Range("N86").Value = -4972.5

# Sheet CRP, Row 89 (Leve Item ID 12584)
$ws_CRP.Range("H89").Value = 2635
$ws_CRP.Range("I89").Value = 2589.25
$ws_CRP.Range("J89").Value = 2726.5
$ws_CRP.Range("K89").Value = 12946.25
$ws_CRP.Range("L89").Value = 13632.5
$ws_CRP.Range("M89").Value = -7330.25
$ws_CRP.Range("N89").Value = -24864.5

# Sheet CUL, Row 51 (Leve Item ID 4646)
$ws_CUL.Range("H51").Value = 2666.5
$ws_CUL.Range("I51").Value = 999
$ws_CUL.Range("J51").Value = 3000
$ws_CUL.Range("K51").Value = 2997
$ws_CUL.Range("L51").Value = 9000
$ws_CUL.Range("M51").Value = -2537
$ws_CUL.Range("N51").Value = -9920

# Sheet CUL, Row 113 (Leve Item ID 27843)
$ws_CUL.Range("H113").Value = 824.5
$ws_CUL.Range("I113").Value = 691.3889
$ws_CUL.Range("J113").Value = 974.25
$ws_CUL.Range("K113").Value = 2074.1667
$ws_CUL.Range("L113").Value = 2922.75
$ws_CUL.Range("M113").Value = 95.83329999999978
$ws_CUL.Range("N113").Value = -7262.75

# Sheet CUL, Row 122 (Leve Item ID 36078)
$ws_CUL.Range("H122").Value = 3476
$ws_CUL.Range("I122").Value = 638.9091
$ws_CUL.Range("J122").Value = 3834.7126
$ws_CUL.Range("K122").Value = 5750.1819
$ws_CUL.Range("L122").Value = 34512.4134
$ws_CUL.Range("M122").Value = -3300.1819
$ws_CUL.Range("N122").Value = -39412.4134

# Sheet CUL, Row 131 (Leve Item ID 36060)
$ws_CUL.Range("H131").Value = 31251498
$ws_CUL.Range("I131").Value = 71430296
$ws_CUL.Range("J131").Value = 1325.3334
$ws_CUL.Range("K131").Value = 214290888
$ws_CUL.Range("L131").Value = 3976.0002
$ws_CUL.Range("M131").Value = -214285848
$ws_CUL.Range("N131").Value = -14056.0002

# Sheet CUL, Row 140 (Leve Item ID 44097)
$ws_CUL.Range("H140").Value = 3425.2942
$ws_CUL.Range("I140").Value = 3702.7273
$ws_CUL.Range("J140").Value = 2916.6667
$ws_CUL.Range("K140").Value = 11108.1819
$ws_CUL.Range("L140").Value = 8750.000100000001
$ws_CUL.Range("M140").Value = -5928.1819
$ws_CUL.Range("N140").Value = -19110.0001

# Sheet GSM, Row 132 (Leve Item ID 44008)
$ws_GSM.Range("H132").Value = 3100.4194
$ws_GSM.Range("I132").Value = 1931.6842
$ws_GSM.Range("J132").Value = 4950.9165
$ws_GSM.Range("K132").Value = 5795.0526
$ws_GSM.Range("L132").Value = 14852.7495
$ws_GSM.Range("M132").Value = -3265.0526
$ws_GSM.Range("N132").Value = -19912.7495

# Sheet WVR, Row 107 (Leve Item ID 27746)
$ws_WVR.Range("H107").Value = 976.6429000000001
$ws_WVR.Range("I107").Value = 669
$ws_WVR.Range("J107").Value = 1745.75
$ws_WVR.Range("K107").Value = 2007
$ws_WVR.Range("L107").Value = 5237.25
$ws_WVR.Range("M107").Value = -87
$ws_WVR.Range("N107").Value = -9077.25

# Sheet WVR, Row 138 (Leve Item ID 42347)
$ws_WVR.Range("H138").Value = 50750
$ws_WVR.Range("J138").Value = 50750
$ws_WVR.Range("L138").Value = 50750
$ws_WVR.Range("N138").Value = -61030

# Sheet WVR, Row 139 (Leve Item ID 43312)
$ws_WVR.Range("H139").Value = 41185.453
$ws_WVR.Range("J139").Value = 41210.953
$ws_WVR.Range("L139").Value = 41210.953
$ws_WVR.Range("N139").Value = -51490.953

# Sheet WVR, Row 140 (Leve Item ID 42506)
$ws_WVR.Range("H140").Value = 42429
$ws_WVR.Range("J140").Value = 42429
$ws_WVR.Range("L140").Value = 42429
$ws_WVR.Range("N140").Value = -52789
